$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P9").Value = 16
$ws.Range("Q9").Value = 6
$ws.Range("AE9").Value = 17

$ws.Range("M11").Value = 22
$ws.Range("N11").Value = 23
$ws.Range("O11").Value = 24
$ws.Range("V11").Value = 40
$ws.Range("W11").Value = 41
$ws.Range("X11").Value = 42
$ws.Range("AN11").Value = 49
$ws.Range("AO11").Value = 50
$ws.Range("AP11").Value = 51
$ws.Range("AW11").Value = 13
$ws.Range("AX11").Value = 14
$ws.Range("AY11").Value = 15

$ws.Range("AJ12:AP12").Clear()

$ws.Range("A9:BC9").Select()
